$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of "Plátano" price data (date serial 44706 = 2022-05-25) is
# inserted right after the current top block, pushing the rest of the
# product's rows down by 4 (one "week" = 4 quality grades here). Insert 4
# blank rows at 1122 (they inherit formatting, incl. the date style on D,
# from the row above) and then populate them.

$ws.Range("A1122:A1125").EntireRow.Insert()

# Row 1122: Maduro
$ws.Cells.Item(1122,1).Value  = 9
$ws.Cells.Item(1122,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1122,3).Value  = "Metropolitana"
$ws.Cells.Item(1122,4).Value  = 44706
$ws.Cells.Item(1122,5).Value  = 13
$ws.Cells.Item(1122,6).Value  = "Fruta"
$ws.Cells.Item(1122,7).Value  = 100108
$ws.Cells.Item(1122,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(1122,9).Value  = 100108006
$ws.Cells.Item(1122,10).Value = "Plátano"
$ws.Cells.Item(1122,11).Value = "Sin especificar"
$ws.Cells.Item(1122,12).Value = "Maduro"
$ws.Cells.Item(1122,13).Value = 880
$ws.Cells.Item(1122,14).Value = 7000
$ws.Cells.Item(1122,15).Value = 7000
$ws.Cells.Item(1122,16).Value = 7000
$ws.Cells.Item(1122,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1122,18).Value = "Ecuador"
$ws.Cells.Item(1122,19).Value = 350
$ws.Cells.Item(1122,20).Value = 20

# Row 1123: Pintón
$ws.Cells.Item(1123,1).Value  = 9
$ws.Cells.Item(1123,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1123,3).Value  = "Metropolitana"
$ws.Cells.Item(1123,4).Value  = 44706
$ws.Cells.Item(1123,5).Value  = 13
$ws.Cells.Item(1123,6).Value  = "Fruta"
$ws.Cells.Item(1123,7).Value  = 100108
$ws.Cells.Item(1123,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(1123,9).Value  = 100108006
$ws.Cells.Item(1123,10).Value = "Plátano"
$ws.Cells.Item(1123,11).Value = "Sin especificar"
$ws.Cells.Item(1123,12).Value = "Pintón"
$ws.Cells.Item(1123,13).Value = 700
$ws.Cells.Item(1123,14).Value = 8000
$ws.Cells.Item(1123,15).Value = 8000
$ws.Cells.Item(1123,16).Value = 8000
$ws.Cells.Item(1123,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1123,18).Value = "Ecuador"
$ws.Cells.Item(1123,19).Value = 400
$ws.Cells.Item(1123,20).Value = 20

# Row 1124: Primera Maduro
$ws.Cells.Item(1124,1).Value  = 9
$ws.Cells.Item(1124,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1124,3).Value  = "Metropolitana"
$ws.Cells.Item(1124,4).Value  = 44706
$ws.Cells.Item(1124,5).Value  = 13
$ws.Cells.Item(1124,6).Value  = "Fruta"
$ws.Cells.Item(1124,7).Value  = 100108
$ws.Cells.Item(1124,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(1124,9).Value  = 100108006
$ws.Cells.Item(1124,10).Value = "Plátano"
$ws.Cells.Item(1124,11).Value = "Sin especificar"
$ws.Cells.Item(1124,12).Value = "Primera Maduro"
$ws.Cells.Item(1124,13).Value = 950
$ws.Cells.Item(1124,14).Value = 9000
$ws.Cells.Item(1124,15).Value = 9000
$ws.Cells.Item(1124,16).Value = 9000
$ws.Cells.Item(1124,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1124,18).Value = "Ecuador"
$ws.Cells.Item(1124,19).Value = 450
$ws.Cells.Item(1124,20).Value = 20

# Row 1125: Primera Pintón
$ws.Cells.Item(1125,1).Value  = 9
$ws.Cells.Item(1125,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1125,3).Value  = "Metropolitana"
$ws.Cells.Item(1125,4).Value  = 44706
$ws.Cells.Item(1125,5).Value  = 13
$ws.Cells.Item(1125,6).Value  = "Fruta"
$ws.Cells.Item(1125,7).Value  = 100108
$ws.Cells.Item(1125,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(1125,9).Value  = 100108006
$ws.Cells.Item(1125,10).Value = "Plátano"
$ws.Cells.Item(1125,11).Value = "Sin especificar"
$ws.Cells.Item(1125,12).Value = "Primera Pintón"
$ws.Cells.Item(1125,13).Value = 750
$ws.Cells.Item(1125,14).Value = 10000
$ws.Cells.Item(1125,15).Value = 10000
$ws.Cells.Item(1125,16).Value = 10000
$ws.Cells.Item(1125,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1125,18).Value = "Ecuador"
$ws.Cells.Item(1125,19).Value = 500
$ws.Cells.Item(1125,20).Value = 20

Write-Host "Inserted 4 new rows at 1122-1125; sheet now spans" $ws.UsedRange.Rows.Count "rows"
